# Auto-generated edit script: updates FFXIV leve-profit market-price snapshots
# across all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), refreshing
# columns H-N (currentAveragePrice*, LevePrice*, LeveProfit*) per the scraped diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1402.579
$ws.Range("I6").Value = 152.85715
$ws.Range("J6").Value = 2131.5833
$ws.Range("K6").Value = 458.57145
$ws.Range("L6").Value = 6394.749899999999
$ws.Range("M6").Value = -346.57145
$ws.Range("N6").Value = -6618.749899999999
$ws.Range("H9").Value = 9091020
$ws.Range("I9").Value = 14285785
$ws.Range("J9").Value = 182.5
$ws.Range("K9").Value = 14285785
$ws.Range("L9").Value = 182.5
$ws.Range("M9").Value = -14285616
$ws.Range("N9").Value = -520.5
$ws.Range("H15").Value = 2168.9167
$ws.Range("I15").Value = 2168.9167
$ws.Range("K15").Value = 6506.750100000001
$ws.Range("M15").Value = -6337.750100000001
$ws.Range("H17").Value = 383003.62
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 383003.62
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 1149010.86
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -1149346.86
$ws.Range("H28").Value = 307.57144
$ws.Range("I28").Value = 277.94446
$ws.Range("J28").Value = 485.33334
$ws.Range("K28").Value = 277.94446
$ws.Range("L28").Value = 485.33334
$ws.Range("M28").Value = 207.05554
$ws.Range("N28").Value = -1455.33334
$ws.Range("H33").Value = 1149.9565
$ws.Range("I33").Value = 678.2143
$ws.Range("J33").Value = 1883.7778
$ws.Range("K33").Value = 678.2143
$ws.Range("L33").Value = 1883.7778
$ws.Range("M33").Value = -449.2143
$ws.Range("N33").Value = -2341.7778
$ws.Range("H53").Value = 302.04544
$ws.Range("I53").Value = 240.75
$ws.Range("K53").Value = 240.75
$ws.Range("M53").Value = 396.25
$ws.Range("H64").Value = 3311.6667
$ws.Range("I64").Value = 4075
$ws.Range("J64").Value = 2930
$ws.Range("K64").Value = 4075
$ws.Range("L64").Value = 2930
$ws.Range("M64").Value = -3827
$ws.Range("N64").Value = -3426
$ws.Range("H67").Value = 3311.6667
$ws.Range("I67").Value = 4075
$ws.Range("J67").Value = 2930
$ws.Range("K67").Value = 4075
$ws.Range("L67").Value = 2930
$ws.Range("M67").Value = -3217
$ws.Range("N67").Value = -4646
$ws.Range("H74").Value = 3894
$ws.Range("I74").Value = 3825
$ws.Range("J74").Value = 3940
$ws.Range("K74").Value = 3825
$ws.Range("L74").Value = 3940
$ws.Range("M74").Value = -2889
$ws.Range("N74").Value = -5812
$ws.Range("H77").Value = 3894
$ws.Range("I77").Value = 3825
$ws.Range("J77").Value = 3940
$ws.Range("K77").Value = 19125
$ws.Range("L77").Value = 19700
$ws.Range("M77").Value = -14445
$ws.Range("N77").Value = -29060
$ws.Range("H100").Value = 4529.091
$ws.Range("I100").Value = 4380
$ws.Range("J100").Value = 4653.3335
$ws.Range("K100").Value = 4380
$ws.Range("L100").Value = 4653.3335
$ws.Range("M100").Value = -3839
$ws.Range("N100").Value = -5735.3335
$ws.Range("H135").Value = 1121.3125
$ws.Range("I135").Value = 649.3077
$ws.Range("K135").Value = 5843.7693
$ws.Range("M135").Value = -3308.7693
$ws.Range("H138").Value = 7068
$ws.Range("I138").Value = 3232.6191
$ws.Range("J138").Value = 8270.135
$ws.Range("K138").Value = 9697.8573
$ws.Range("L138").Value = 24810.405
$ws.Range("M138").Value = -4557.8573
$ws.Range("N138").Value = -35090.405

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2035.5491
$ws.Range("I132").Value = 1558.359
$ws.Range("J132").Value = 3586.4167
$ws.Range("K132").Value = 4675.076999999999
$ws.Range("L132").Value = 10759.2501
$ws.Range("M132").Value = -2145.076999999999
$ws.Range("N132").Value = -15819.2501

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 2500
$ws.Range("J15").Value = 2500
$ws.Range("L15").Value = 2500
$ws.Range("N15").Value = -2954
$ws.Range("H134").Value = 2899.5217
$ws.Range("I134").Value = 2796.4146
$ws.Range("J134").Value = 3745
$ws.Range("K134").Value = 8389.2438
$ws.Range("L134").Value = 11235
$ws.Range("M134").Value = -5854.2438
$ws.Range("N134").Value = -16305

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 10002850
$ws.Range("I58").Value = 1430.6765
$ws.Range("J58").Value = 31255868
$ws.Range("K58").Value = 1430.6765
$ws.Range("L58").Value = 31255868
$ws.Range("M58").Value = -1227.6765
$ws.Range("N58").Value = -31256274
$ws.Range("H99").Value = 4035.5
$ws.Range("I99").Value = 1233
$ws.Range("K99").Value = 1233
$ws.Range("M99").Value = 265
$ws.Range("H126").Value = 4035.5
$ws.Range("I126").Value = 1233
$ws.Range("K126").Value = 3699
$ws.Range("M126").Value = -1229
$ws.Range("H132").Value = 3548.8262
$ws.Range("I132").Value = 2634.0667
$ws.Range("J132").Value = 5264
$ws.Range("K132").Value = 7902.2001
$ws.Range("L132").Value = 15792
$ws.Range("M132").Value = -5372.2001
$ws.Range("N132").Value = -20852
$ws.Range("H134").Value = 2895.1353
$ws.Range("I134").Value = 2310.4546
$ws.Range("J134").Value = 3752.6667
$ws.Range("K134").Value = 6931.3638
$ws.Range("L134").Value = 11258.0001
$ws.Range("M134").Value = -4396.3638
$ws.Range("N134").Value = -16328.0001
$ws.Range("H136").Value = 10002850
$ws.Range("I136").Value = 1430.6765
$ws.Range("J136").Value = 31255868
$ws.Range("K136").Value = 4292.029500000001
$ws.Range("L136").Value = 93767604
$ws.Range("M136").Value = -1742.029500000001
$ws.Range("N136").Value = -93772704

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 125
$ws.Range("I7").Value = 50
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 150
$ws.Range("L7").Value = 600
$ws.Range("M7").Value = -38
$ws.Range("N7").Value = -824
$ws.Range("H33").Value = 114.75
$ws.Range("I33").Value = 54
$ws.Range("J33").Value = 175.5
$ws.Range("K33").Value = 324
$ws.Range("L33").Value = 1053
$ws.Range("M33").Value = -41
$ws.Range("N33").Value = -1619
$ws.Range("H34").Value = 6105.421
$ws.Range("I34").Value = 161.25
$ws.Range("J34").Value = 10428.454
$ws.Range("K34").Value = 483.75
$ws.Range("L34").Value = 31285.362
$ws.Range("M34").Value = -399.75
$ws.Range("N34").Value = -31453.362
$ws.Range("H68").Value = 649.7778
$ws.Range("I68").Value = 769.6
$ws.Range("J68").Value = 500
$ws.Range("K68").Value = 2308.8
$ws.Range("L68").Value = 1500
$ws.Range("M68").Value = -1497.8
$ws.Range("N68").Value = -3122
$ws.Range("H71").Value = 649.7778
$ws.Range("I71").Value = 769.6
$ws.Range("J71").Value = 500
$ws.Range("K71").Value = 6926.400000000001
$ws.Range("L71").Value = 4500
$ws.Range("M71").Value = -2870.400000000001
$ws.Range("N71").Value = -12612
$ws.Range("H80").Value = 2736.7273
$ws.Range("I80").Value = 2034.6666
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 6103.9998
$ws.Range("L80").Value = 9000
$ws.Range("M80").Value = -5167.9998
$ws.Range("N80").Value = -10872
$ws.Range("H83").Value = 2736.7273
$ws.Range("I83").Value = 2034.6666
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 18311.9994
$ws.Range("L83").Value = 27000
$ws.Range("M83").Value = -13631.9994
$ws.Range("N83").Value = -36360
$ws.Range("H92").Value = 1600.4546
$ws.Range("J92").Value = 1800.3334
$ws.Range("L92").Value = 5401.0002
$ws.Range("N92").Value = -7897.0002
$ws.Range("H109").Value = 727.0769
$ws.Range("I109").Value = 426.1
$ws.Range("J109").Value = 1730.3334
$ws.Range("K109").Value = 1278.3
$ws.Range("L109").Value = 5191.0002
$ws.Range("M109").Value = -238.3000000000002
$ws.Range("N109").Value = -7271.0002
$ws.Range("H113").Value = 3572463.8
$ws.Range("I113").Value = 10000696
$ws.Range("J113").Value = 1223.8334
$ws.Range("K113").Value = 30002088
$ws.Range("L113").Value = 3671.5002
$ws.Range("M113").Value = -29999918
$ws.Range("N113").Value = -8011.5002
$ws.Range("H131").Value = 1708.0233
$ws.Range("J131").Value = 1266.125
$ws.Range("L131").Value = 3798.375
$ws.Range("N131").Value = -13878.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1668.6364
$ws.Range("I113").Value = 794.375
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 794.375
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = 1375.625
$ws.Range("N113").Value = -8340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4000
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H40").Value = 16666.334
$ws.Range("I40").Value = 22499.5
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 22499.5
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -22363.5
$ws.Range("N40").Value = -5272
$ws.Range("H126").Value = 4000
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 794
$ws.Range("I81").Value = 780.4
$ws.Range("J81").Value = 816.6667
$ws.Range("K81").Value = 1560.8
$ws.Range("L81").Value = 1633.3334
$ws.Range("M81").Value = -499.8
$ws.Range("N81").Value = -3755.3334
$ws.Range("H84").Value = 794
$ws.Range("I84").Value = 780.4
$ws.Range("J84").Value = 816.6667
$ws.Range("K84").Value = 7804
$ws.Range("L84").Value = 8166.666999999999
$ws.Range("M84").Value = -2500
$ws.Range("N84").Value = -18774.667
$ws.Range("H126").Value = 60744.293
$ws.Range("I126").Value = 77436.38
$ws.Range("J126").Value = 6495
$ws.Range("K126").Value = 232309.14
$ws.Range("L126").Value = 19485
$ws.Range("M126").Value = -229839.14
$ws.Range("N126").Value = -24425
$ws.Range("H132").Value = 2236013.8
$ws.Range("I132").Value = 3035564
$ws.Range("J132").Value = 37250.168
$ws.Range("K132").Value = 9106692
$ws.Range("L132").Value = 111750.504
$ws.Range("M132").Value = -9104162
$ws.Range("N132").Value = -116810.504
$ws.Range("H140").Value = 72557.25
$ws.Range("J140").Value = 72557.25
$ws.Range("L140").Value = 72557.25
$ws.Range("N140").Value = -82917.25
$ws.Range("H141").Value = 28200
$ws.Range("J141").Value = 28200
$ws.Range("L141").Value = 28200
$ws.Range("N141").Value = -38560
